$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update D,E,F,G,H (open_price/close_price/high_price/low_price/shares_outstanding)
# and I (fixed_ticker) for rows 2-44 to reflect corrected WDAY price history

# Row 2
$ws.Cells.Item(2, 4).Value = 79.65000152587891
$ws.Cells.Item(2, 5).Value = 85.5
$ws.Cells.Item(2, 6).Value = 95.16999816894533
$ws.Cells.Item(2, 7).Value = 78.40000152587891
$ws.Cells.Item(2, 8).Value = 217000000
$ws.Cells.Item(2, 9).Value = "WDAY"

# Row 3
$ws.Cells.Item(3, 4).Value = 91.20999908447266
$ws.Cells.Item(3, 5).Value = 78.91999816894531
$ws.Cells.Item(3, 6).Value = 93.62000274658205
$ws.Cells.Item(3, 7).Value = 78.45999908447266
$ws.Cells.Item(3, 8).Value = 217000000
$ws.Cells.Item(3, 9).Value = "WDAY"

# Row 4
$ws.Cells.Item(4, 4).Value = 84.26999664306641
$ws.Cells.Item(4, 5).Value = 70.26000213623047
$ws.Cells.Item(4, 6).Value = 84.83999633789062
$ws.Cells.Item(4, 7).Value = 65.33000183105469
$ws.Cells.Item(4, 8).Value = 217000000
$ws.Cells.Item(4, 9).Value = "WDAY"

# Row 5
$ws.Cells.Item(5, 4).Value = 79.48000335693359
$ws.Cells.Item(5, 5).Value = 83.70999908447266
$ws.Cells.Item(5, 6).Value = 85.66999816894531
$ws.Cells.Item(5, 7).Value = 77.90000152587891
$ws.Cells.Item(5, 8).Value = 217000000
$ws.Cells.Item(5, 9).Value = "WDAY"

# Row 6
$ws.Cells.Item(6, 4).Value = 62.56000137329102
$ws.Cells.Item(6, 5).Value = 60.45000076293945
$ws.Cells.Item(6, 6).Value = 65.90000152587891
$ws.Cells.Item(6, 7).Value = 47.31999969482422
$ws.Cells.Item(6, 8).Value = 217000000
$ws.Cells.Item(6, 9).Value = "WDAY"

# Row 7
$ws.Cells.Item(7, 4).Value = 74.94999694824219
$ws.Cells.Item(7, 5).Value = 75.83999633789062
$ws.Cells.Item(7, 6).Value = 78.93000030517578
$ws.Cells.Item(7, 7).Value = 69
$ws.Cells.Item(7, 8).Value = 217000000
$ws.Cells.Item(7, 9).Value = "WDAY"

# Row 8
$ws.Cells.Item(8, 4).Value = 83.33999633789062
$ws.Cells.Item(8, 5).Value = 84.79000091552734
$ws.Cells.Item(8, 6).Value = 85.95999908447266
$ws.Cells.Item(8, 7).Value = 78.33999633789062
$ws.Cells.Item(8, 8).Value = 217000000
$ws.Cells.Item(8, 9).Value = "WDAY"

# Row 9
$ws.Cells.Item(9, 4).Value = 87.22000122070312
$ws.Cells.Item(9, 5).Value = 84.31999969482422
$ws.Cells.Item(9, 6).Value = 87.26000213623047
$ws.Cells.Item(9, 7).Value = 78.05000305175781
$ws.Cells.Item(9, 8).Value = 217000000
$ws.Cells.Item(9, 9).Value = "WDAY"

# Row 10
$ws.Cells.Item(10, 4).Value = 83.91999816894531
$ws.Cells.Item(10, 5).Value = 82.93000030517578
$ws.Cells.Item(10, 6).Value = 92.45999908447266
$ws.Cells.Item(10, 7).Value = 81.45999908447266
$ws.Cells.Item(10, 8).Value = 217000000
$ws.Cells.Item(10, 9).Value = "WDAY"

# Row 11
$ws.Cells.Item(11, 4).Value = 87.62000274658203
$ws.Cells.Item(11, 5).Value = 99.98000335693359
$ws.Cells.Item(11, 6).Value = 101.3300018310547
$ws.Cells.Item(11, 7).Value = 87.08999633789062
$ws.Cells.Item(11, 8).Value = 217000000
$ws.Cells.Item(11, 9).Value = "WDAY"

# Row 12
$ws.Cells.Item(12, 4).Value = 103.129997253418
$ws.Cells.Item(12, 5).Value = 109.6900024414062
$ws.Cells.Item(12, 6).Value = 111.4499969482422
$ws.Cells.Item(12, 7).Value = 96.83000183105467
$ws.Cells.Item(12, 8).Value = 217000000
$ws.Cells.Item(12, 9).Value = "WDAY"

# Row 13
$ws.Cells.Item(13, 4).Value = 111.9400024414062
$ws.Cells.Item(13, 5).Value = 103
$ws.Cells.Item(13, 6).Value = 116.8899993896484
$ws.Cells.Item(13, 7).Value = 102.1100006103516
$ws.Cells.Item(13, 8).Value = 217000000
$ws.Cells.Item(13, 9).Value = "WDAY"

# Row 14
$ws.Cells.Item(14, 4).Value = 119.0299987792969
$ws.Cells.Item(14, 5).Value = 126.6699981689453
$ws.Cells.Item(14, 6).Value = 131.1799926757812
$ws.Cells.Item(14, 7).Value = 107.75
$ws.Cells.Item(14, 8).Value = 217000000
$ws.Cells.Item(14, 9).Value = "WDAY"

# Row 15
$ws.Cells.Item(15, 4).Value = 123.9700012207031
$ws.Cells.Item(15, 5).Value = 130.9600067138672
$ws.Cells.Item(15, 6).Value = 138.6499938964844
$ws.Cells.Item(15, 7).Value = 123.1500015258789
$ws.Cells.Item(15, 8).Value = 217000000
$ws.Cells.Item(15, 9).Value = "WDAY"

# Row 16
$ws.Cells.Item(16, 4).Value = 123.6100006103516
$ws.Cells.Item(16, 5).Value = 154.5399932861328
$ws.Cells.Item(16, 6).Value = 155.8099975585938
$ws.Cells.Item(16, 7).Value = 122.9199981689453
$ws.Cells.Item(16, 8).Value = 217000000
$ws.Cells.Item(16, 9).Value = "WDAY"

# Row 17
$ws.Cells.Item(17, 4).Value = 133.5500030517578
$ws.Cells.Item(17, 5).Value = 164
$ws.Cells.Item(17, 6).Value = 166.3600006103516
$ws.Cells.Item(17, 7).Value = 117.7200012207031
$ws.Cells.Item(17, 8).Value = 217000000
$ws.Cells.Item(17, 9).Value = "WDAY"

# Row 18
$ws.Cells.Item(18, 4).Value = 180.6699981689453
$ws.Cells.Item(18, 5).Value = 197.9299926757812
$ws.Cells.Item(18, 6).Value = 199.4499969482422
$ws.Cells.Item(18, 7).Value = 180.3300018310547
$ws.Cells.Item(18, 8).Value = 217000000
$ws.Cells.Item(18, 9).Value = "WDAY"

# Row 19
$ws.Cells.Item(19, 4).Value = 206.5
$ws.Cells.Item(19, 5).Value = 204.1199951171875
$ws.Cells.Item(19, 6).Value = 217.6300048828125
$ws.Cells.Item(19, 7).Value = 191.25
$ws.Cells.Item(19, 8).Value = 217000000
$ws.Cells.Item(19, 9).Value = "WDAY"

# Row 20
$ws.Cells.Item(20, 4).Value = 201
$ws.Cells.Item(20, 5).Value = 177.2799987792969
$ws.Cells.Item(20, 6).Value = 206.759994506836
$ws.Cells.Item(20, 7).Value = 173.1999969482422
$ws.Cells.Item(20, 8).Value = 217000000
$ws.Cells.Item(20, 9).Value = "WDAY"

# Row 21
$ws.Cells.Item(21, 4).Value = 163.6300048828125
$ws.Cells.Item(21, 5).Value = 179.1199951171875
$ws.Cells.Item(21, 6).Value = 180.0299987792969
$ws.Cells.Item(21, 7).Value = 158.7899932861328
$ws.Cells.Item(21, 8).Value = 217000000
$ws.Cells.Item(21, 9).Value = "WDAY"

# Row 22
$ws.Cells.Item(22, 4).Value = 185.6000061035156
$ws.Cells.Item(22, 5).Value = 173.25
$ws.Cells.Item(22, 6).Value = 202
$ws.Cells.Item(22, 7).Value = 164.6999969482422
$ws.Cells.Item(22, 8).Value = 217000000
$ws.Cells.Item(22, 9).Value = "WDAY"

# Row 23
$ws.Cells.Item(23, 4).Value = 150.4400024414062
$ws.Cells.Item(23, 5).Value = 183.4299926757812
$ws.Cells.Item(23, 6).Value = 187
$ws.Cells.Item(23, 7).Value = 144.8099975585938
$ws.Cells.Item(23, 8).Value = 217000000
$ws.Cells.Item(23, 9).Value = "WDAY"

# Row 24
$ws.Cells.Item(24, 4).Value = 182.4700012207031
$ws.Cells.Item(24, 5).Value = 239.7100067138672
$ws.Cells.Item(24, 6).Value = 248.75
$ws.Cells.Item(24, 7).Value = 174.5200042724609
$ws.Cells.Item(24, 8).Value = 217000000
$ws.Cells.Item(24, 9).Value = "WDAY"

# Row 25
$ws.Cells.Item(25, 4).Value = 209.729995727539
$ws.Cells.Item(25, 5).Value = 224.7899932861328
$ws.Cells.Item(25, 6).Value = 231.9400024414062
$ws.Cells.Item(25, 7).Value = 201.6199951171875
$ws.Cells.Item(25, 8).Value = 217000000
$ws.Cells.Item(25, 9).Value = "WDAY"

# Row 26
$ws.Cells.Item(26, 4).Value = 228.5200042724609
$ws.Cells.Item(26, 5).Value = 245.1799926757812
$ws.Cells.Item(26, 6).Value = 282.7699890136719
$ws.Cells.Item(26, 7).Value = 228.5200042724609
$ws.Cells.Item(26, 8).Value = 217000000
$ws.Cells.Item(26, 9).Value = "WDAY"

# Row 27
$ws.Cells.Item(27, 4).Value = 248.3500061035156
$ws.Cells.Item(27, 5).Value = 228.7200012207031
$ws.Cells.Item(27, 6).Value = 250.1900024414062
$ws.Cells.Item(27, 7).Value = 217.6000061035156
$ws.Cells.Item(27, 8).Value = 217000000
$ws.Cells.Item(27, 9).Value = "WDAY"

# Row 28
$ws.Cells.Item(28, 4).Value = 235.1000061035156
$ws.Cells.Item(28, 5).Value = 273.1600036621094
$ws.Cells.Item(28, 6).Value = 275.7099914550781
$ws.Cells.Item(28, 7).Value = 228.0700073242188
$ws.Cells.Item(28, 8).Value = 217000000
$ws.Cells.Item(28, 9).Value = "WDAY"

# Row 29
$ws.Cells.Item(29, 4).Value = 291.1099853515625
$ws.Cells.Item(29, 5).Value = 274.2300109863281
$ws.Cells.Item(29, 6).Value = 307.8099975585937
$ws.Cells.Item(29, 7).Value = 268.3699951171875
$ws.Cells.Item(29, 8).Value = 217000000
$ws.Cells.Item(29, 9).Value = "WDAY"

# Row 30
$ws.Cells.Item(30, 4).Value = 256.9599914550781
$ws.Cells.Item(30, 5).Value = 229.0500030517578
$ws.Cells.Item(30, 6).Value = 257
$ws.Cells.Item(30, 7).Value = 205.8999938964844
$ws.Cells.Item(30, 8).Value = 217000000
$ws.Cells.Item(30, 9).Value = "WDAY"

# Row 31
$ws.Cells.Item(31, 4).Value = 205.7799987792969
$ws.Cells.Item(31, 5).Value = 156.3000030517578
$ws.Cells.Item(31, 6).Value = 209.3699951171875
$ws.Cells.Item(31, 7).Value = 149.0599975585938
$ws.Cells.Item(31, 8).Value = 217000000
$ws.Cells.Item(31, 9).Value = "WDAY"

# Row 32
$ws.Cells.Item(32, 4).Value = 152.9700012207031
$ws.Cells.Item(32, 5).Value = 164.5599975585938
$ws.Cells.Item(32, 6).Value = 179.8300018310547
$ws.Cells.Item(32, 7).Value = 151.6999969482422
$ws.Cells.Item(32, 8).Value = 217000000
$ws.Cells.Item(32, 9).Value = "WDAY"

# Row 33
$ws.Cells.Item(33, 4).Value = 158.4799957275391
$ws.Cells.Item(33, 5).Value = 167.8999938964844
$ws.Cells.Item(33, 6).Value = 167.9299926757812
$ws.Cells.Item(33, 7).Value = 128.7200012207031
$ws.Cells.Item(33, 8).Value = 217000000
$ws.Cells.Item(33, 9).Value = "WDAY"

# Row 34
$ws.Cells.Item(34, 4).Value = 180.5800018310547
$ws.Cells.Item(34, 5).Value = 185.4700012207031
$ws.Cells.Item(34, 6).Value = 193.6399993896484
$ws.Cells.Item(34, 7).Value = 177.7799987792969
$ws.Cells.Item(34, 8).Value = 217000000
$ws.Cells.Item(34, 9).Value = "WDAY"

# Row 35
$ws.Cells.Item(35, 4).Value = 185.509994506836
$ws.Cells.Item(35, 5).Value = 211.9900054931641
$ws.Cells.Item(35, 6).Value = 218.8800048828125
$ws.Cells.Item(35, 7).Value = 174.25
$ws.Cells.Item(35, 8).Value = 217000000
$ws.Cells.Item(35, 9).Value = "WDAY"

# Row 36
$ws.Cells.Item(36, 4).Value = 235.8399963378907
$ws.Cells.Item(36, 5).Value = 244.5
$ws.Cells.Item(36, 6).Value = 245.6000061035156
$ws.Cells.Item(36, 7).Value = 218.4400024414062
$ws.Cells.Item(36, 8).Value = 217000000
$ws.Cells.Item(36, 9).Value = "WDAY"

# Row 37
$ws.Cells.Item(37, 4).Value = 209.4600067138672
$ws.Cells.Item(37, 5).Value = 270.7200012207031
$ws.Cells.Item(37, 6).Value = 271.989990234375
$ws.Cells.Item(37, 7).Value = 206.9900054931641
$ws.Cells.Item(37, 8).Value = 217000000
$ws.Cells.Item(37, 9).Value = "WDAY"

# Row 38
$ws.Cells.Item(38, 4).Value = 292.7000122070312
$ws.Cells.Item(38, 5).Value = 294.6600036621094
$ws.Cells.Item(38, 6).Value = 311.2799987792969
$ws.Cells.Item(38, 7).Value = 284.6300048828125
$ws.Cells.Item(38, 8).Value = 217000000
$ws.Cells.Item(38, 9).Value = "WDAY"

# Row 39
$ws.Cells.Item(39, 4).Value = 245.3800048828125
$ws.Cells.Item(39, 5).Value = 211.4499969482422
$ws.Cells.Item(39, 6).Value = 263.6600036621094
$ws.Cells.Item(39, 7).Value = 204.5200042724609
$ws.Cells.Item(39, 8).Value = 217000000
$ws.Cells.Item(39, 9).Value = "WDAY"

# Row 40
$ws.Cells.Item(40, 4).Value = 227.3300018310547
$ws.Cells.Item(40, 5).Value = 263.1900024414062
$ws.Cells.Item(40, 6).Value = 266.8200073242188
$ws.Cells.Item(40, 7).Value = 199.8099975585937
$ws.Cells.Item(40, 8).Value = 217000000
$ws.Cells.Item(40, 9).Value = "WDAY"

# Row 41
$ws.Cells.Item(41, 4).Value = 233.1699981689453
$ws.Cells.Item(41, 5).Value = 249.9900054931641
$ws.Cells.Item(41, 6).Value = 278.9800109863281
$ws.Cells.Item(41, 7).Value = 231.7400054931641
$ws.Cells.Item(41, 8).Value = 217000000
$ws.Cells.Item(41, 9).Value = "WDAY"

# Row 42
$ws.Cells.Item(42, 4).Value = 259.7300109863281
$ws.Cells.Item(42, 5).Value = 263.3399963378906
$ws.Cells.Item(42, 6).Value = 283.6799926757812
$ws.Cells.Item(42, 7).Value = 251.0299987792969
$ws.Cells.Item(42, 8).Value = 217000000
$ws.Cells.Item(42, 9).Value = "WDAY"

# Row 43
$ws.Cells.Item(43, 4).Value = 247.0500030517578
$ws.Cells.Item(43, 5).Value = 247.7100067138672
$ws.Cells.Item(43, 6).Value = 276
$ws.Cells.Item(43, 7).Value = 237
$ws.Cells.Item(43, 8).Value = 217000000
$ws.Cells.Item(43, 9).Value = "WDAY"

# Row 44
$ws.Cells.Item(44, 4).Value = 228.7400054931641
$ws.Cells.Item(44, 5).Value = 230.8200073242188
$ws.Cells.Item(44, 6).Value = 234.9199981689453
$ws.Cells.Item(44, 7).Value = 206.7700042724609
$ws.Cells.Item(44, 8).Value = 217000000
$ws.Cells.Item(44, 9).Value = "WDAY"
